$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Language Code" header column (C1), matching the style used by
# the existing header cells A1/B1 (bold font, grey fill, border).
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Language Code"

# Match the column width used by the other header columns.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Move/collapse the selection to A2, as in the saved workbook.
$ws.Range("A2").Select()
